# Generate Report for Handback
# Refresh the handoff/handback timestamp values that get written when the
# handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 3e3e0a4d... row
$wsOverview.Range("G2").Value = "2016-08-30 09:14:34"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 3e3e0a4d... row
$wsZhCn.Range("H2").Value = "2016-08-30 09:14:23"
$wsZhCn.Range("K2").Value = "2016-08-30 09:15:16"

# de-de sheet: "Correspond Handoff Datetime" (mirrors the Overview value)
# and "Correspond Handback DateTime" for the 3e3e0a4d... row
$wsDeDe.Range("H2").Value = "2016-08-30 09:14:34"
$wsDeDe.Range("K2").Value = "2016-08-30 09:15:36"
